# "Multiple fd intrest payment (Rework)"
#
# This adds a new worksheet "Bondissue" right after the existing "BondApp"
# sheet (making it the new last sheet / newly active tab), seeded with the
# same header/data row layout as "BondApp" plus one extra "Reject" column,
# and a couple of tweaked data values ("Bond_Issue" / "A" instead of
# "Bond_Application" / "AAA").

$wb  = $excel.ActiveWorkbook
$ws20 = $wb.Worksheets.Item("BondApp")

# --- Add the new sheet right after BondApp -------------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws20)
$newSheet.Name = "Bondissue"

# --- Seed it with BondApp's header + data row (values, then formats) ----
$ws20.Range("A1:O2").Copy()
$newSheet.Range("A1").PasteSpecial(-4163)   # xlPasteValuesAndNumberFormats
$ws20.Range("A1:O2").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)   # xlPasteFormats

$newSheet.Rows.Item(1).RowHeight = 45
$newSheet.Rows.Item(2).RowHeight = 30

# --- Tweak a couple of data cells for the new sheet ----------------------
$newSheet.Range("A2").Value2 = "Bond_Issue"
$newSheet.Range("E2").Value2 = "A"

# --- Add the new "Reject" column (P) -------------------------------------
$newSheet.Range("P1").Value2 = "Reject"
$newSheet.Range("P2").Value2 = "Rejected"
$ws20.Range("O1").Copy()
$newSheet.Range("P1").PasteSpecial(-4122)   # match header formatting

# --- Fix up sheet view / selection state ---------------------------------
# BondApp is no longer the active/selected tab; its selection now just
# spans the used range instead of the old "M8" cell.
[void]$ws20.Range("A1:O2").Select()

# The new sheet becomes the active tab, with the cursor just past the
# last used column.
[void]$newSheet.Range("Q1").Select()
